$d = $word.ActiveDocument

$replacements = @(
    @("2023-12-13 Wednesday", "2023-12-14 Thursday"),
    @("42×49=", "91×76="),
    @("63×98=", "31×61="),
    @("57×19=", "56×46="),
    @("15×41=", "62×83="),
    @("80×92=", "93×12="),
    @("12×31=", "46×39="),
    @("98×64=", "54×26="),
    @("36×59=", "65×77="),
    @("72×26=", "80×21="),
    @("61×33=", "93×15="),
    @("79×76=", "36×76="),
    @("86×49=", "65×76="),
    @("54×73=", "13×91="),
    @("25×68=", "76×62="),
    @("11×94=", "30×99="),
    @("86×29=", "55×97="),
    @("87×22=", "47×21="),
    @("96×47=", "12×99="),
    @("66×58=", "36×22="),
    @("51×63=", "65×32="),
    @("71×22=", "91×62="),
    @("71×79=", "53×67="),
    @("43×53=", "83×55="),
    @("19×37=", "74×90="),
    @("12×87=", "55×56=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
